$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update filenames in columns D and E (rows 2-4)
$ws.Range("D2:D4").Value = "TC12_CDS_Filter_InstrumentModel-NextSeq 500_Neo4jData.xlsx"
$ws.Range("E2:E4").Value = "TC12_CDS_Filter_InstrumentModel-NextSeq 500_WebData.xlsx"

# Update Cypher queries in column B (rows 2-4) - instrument model filter text
$ws.Range("B2").Value = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['NextSeq 500']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, collect(distinct samp.sample_id) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY ``Participant ID``LIMIT 100"
$ws.Range("B3").Value = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['NextSeq 500']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"
$ws.Range("B4").Value = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['NextSeq 500']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN `n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER By f.file_name LIMIT 100"

# Update Cypher queries in column C (rows 2-4) - same summary query text
$ws.Range("C2:C4").Value = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['NextSeq 500']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

# Column D width adjustment (autofit-driven widen due to longer content)
$ws.Columns("D").ColumnWidth = 84.6

# Update selected cell to D4
$ws.Range("D4").Select()
